# update scripts wuth new tpm
#
# The NATMI ligand-receptor TPM recompute changed which "sending cluster"
# rows exist for Il11 -> Il11ra1 (only MuSCs now sends, the old ECs rows are
# gone) and refreshed the specificity / expression numbers for the
# remaining MuSCs rows (one row per target cluster: ECs, FAPs, MuSCs).
#
# Before: 6 data rows (rows 2-7) -> ECs x {ECs,FAPs,MuSCs}, MuSCs x {ECs,FAPs,MuSCs}
# After : 3 data rows (rows 2-4) -> MuSCs x {ECs,FAPs,MuSCs}, with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old trailing "MuSCs" block (rows 5-7); the surviving rows 2-4 get
# overwritten below with the refreshed MuSCs data, and the old "ECs" block
# disappears entirely.
$ws.Rows("5:7").Delete()

# Row 2: MuSCs -> Il11 -> Il11ra1 -> ECs
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Il11"
$ws.Range("C2").Value = "Il11ra1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.15518
$ws.Range("H2").Value = 0.46554
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.113196333333334
$ws.Range("N2").Value = 12.339589
$ws.Range("O2").Value = 0.06757131683644607
$ws.Range("P2").Value = 0.06757131683644607
$ws.Range("Q2").Value = 0.6382858070066668
$ws.Range("R2").Value = 5.74457226306
$ws.Range("S2").Value = 0.06757131683644607
$ws.Range("T2").Value = 0.06757131683644607

# Row 3: MuSCs -> Il11 -> Il11ra1 -> FAPs
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Il11"
$ws.Range("C3").Value = "Il11ra1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.15518
$ws.Range("H3").Value = 0.46554
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 48.58222466666667
$ws.Range("N3").Value = 145.746674
$ws.Range("O3").Value = 0.7981055679173932
$ws.Range("P3").Value = 0.7981055679173932
$ws.Range("Q3").Value = 7.538989623773334
$ws.Range("R3").Value = 67.85090661396
$ws.Range("S3").Value = 0.7981055679173932
$ws.Range("T3").Value = 0.7981055679173932

# Row 4: MuSCs -> Il11 -> Il11ra1 -> MuSCs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Il11"
$ws.Range("C4").Value = "Il11ra1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.15518
$ws.Range("H4").Value = 0.46554
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8.176507
$ws.Range("N4").Value = 24.529521
$ws.Range("O4").Value = 0.1343231152461608
$ws.Range("P4").Value = 0.1343231152461608
$ws.Range("Q4").Value = 1.26883035626
$ws.Range("R4").Value = 11.41947320634
$ws.Range("S4").Value = 0.1343231152461608
$ws.Range("T4").Value = 0.1343231152461608
